# Commit: "case with 380 kV done"
# Updates Code/Results/Cases/Case_5_243/res_bus/vm_pu.xlsx (Sheet1).
#
# The slack-bus setpoint (column B) moves from 1.05 p.u. to 1.02 p.u. (380 kV case),
# and the resulting bus voltage magnitudes for all other buses (columns C-F, I-N)
# are replaced with the freshly recomputed power-flow results for rows 2-25.
# Column A (row index), column G (fixed at 1 p.u.) and the header row are unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.03055029072378
$ws.Cells.Item(2, 4).Value = 1.039582631608319
$ws.Cells.Item(2, 5).Value = 1.030233735645875
$ws.Cells.Item(2, 6).Value = 1.048059835917407
$ws.Cells.Item(2, 9).Value = 1.034988157074656
$ws.Cells.Item(2, 10).Value = 1.035690915942963
$ws.Cells.Item(2, 11).Value = 1.04236740876327
$ws.Cells.Item(2, 12).Value = 1.033045339321515
$ws.Cells.Item(2, 13).Value = 1.050820736706897
$ws.Cells.Item(2, 14).Value = 1.005712725503983

# Row 3
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.031525713853264
$ws.Cells.Item(3, 4).Value = 1.040334047593937
$ws.Cells.Item(3, 5).Value = 1.031063231416374
$ws.Cells.Item(3, 6).Value = 1.04894409824415
$ws.Cells.Item(3, 9).Value = 1.035156002880945
$ws.Cells.Item(3, 10).Value = 1.036307707967146
$ws.Cells.Item(3, 11).Value = 1.04292922803617
$ws.Cells.Item(3, 12).Value = 1.033683090194355
$ws.Cells.Item(3, 13).Value = 1.05151678109444

# Row 4
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.032157419879803
$ws.Cells.Item(4, 4).Value = 1.040820539874346
$ws.Cells.Item(4, 5).Value = 1.031600807132731
$ws.Cells.Item(4, 6).Value = 1.049516877315329
$ws.Cells.Item(4, 9).Value = 1.035263397277794
$ws.Cells.Item(4, 10).Value = 1.036706761791015
$ws.Cells.Item(4, 11).Value = 1.043292374736066
$ws.Cells.Item(4, 12).Value = 1.034095948083048
$ws.Cells.Item(4, 13).Value = 1.051967133040545

# Row 5
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.032423117985251
$ws.Cells.Item(5, 4).Value = 1.041025125892819
$ws.Cells.Item(5, 5).Value = 1.031827002915325
$ws.Cells.Item(5, 6).Value = 1.049757816376128
$ws.Cells.Item(5, 9).Value = 1.035308254932027
$ws.Cells.Item(5, 10).Value = 1.036874510650811
$ws.Cells.Item(5, 11).Value = 1.043444947683672
$ws.Cells.Item(5, 12).Value = 1.03426955811229
$ws.Cells.Item(5, 13).Value = 1.052156451330734

# Row 6
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.032467737399309
$ws.Cells.Item(6, 4).Value = 1.04105948053718
$ws.Cells.Item(6, 5).Value = 1.031864993811216
$ws.Cells.Item(6, 6).Value = 1.049798279433289
$ws.Cells.Item(6, 9).Value = 1.035315769652412
$ws.Cells.Item(6, 10).Value = 1.03690267561095
$ws.Cells.Item(6, 11).Value = 1.043470559814552
$ws.Cells.Item(6, 12).Value = 1.034298710584372
$ws.Cells.Item(6, 13).Value = 1.05218823810251

# Row 7
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.032160969645518
$ws.Cells.Item(7, 4).Value = 1.040823273310394
$ws.Cells.Item(7, 5).Value = 1.031603828791743
$ws.Cells.Item(7, 6).Value = 1.049520096195782
$ws.Cells.Item(7, 9).Value = 1.035263997811875
$ws.Cells.Item(7, 10).Value = 1.036709003313209
$ws.Cells.Item(7, 11).Value = 1.043294413793641
$ws.Cells.Item(7, 12).Value = 1.034098267694889
$ws.Cells.Item(7, 13).Value = 1.051969662760316

# Row 8
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.030879826778748
$ws.Cells.Item(8, 4).Value = 1.039836518200812
$ws.Cells.Item(8, 5).Value = 1.030513893640531
$ws.Cells.Item(8, 6).Value = 1.048358550956346
$ws.Cells.Item(8, 9).Value = 1.035045132198321
$ws.Cells.Item(8, 10).Value = 1.035899373860225
$ws.Cells.Item(8, 11).Value = 1.042557357909368
$ws.Cells.Item(8, 12).Value = 1.033260830174246
$ws.Cells.Item(8, 13).Value = 1.051055974749293

# Row 9
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.028626479261632
$ws.Cells.Item(9, 4).Value = 1.03809990375174
$ws.Cells.Item(9, 5).Value = 1.028599752738653
$ws.Cells.Item(9, 6).Value = 1.046316435341159
$ws.Cells.Item(9, 9).Value = 1.034650196283176
$ws.Cells.Item(9, 10).Value = 1.034472348170658
$ws.Cells.Item(9, 11).Value = 1.041255646555638
$ws.Cells.Item(9, 12).Value = 1.031786661570097
$ws.Cells.Item(9, 13).Value = 1.049445722806262

# Row 10
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.027127115180575
$ws.Cells.Item(10, 4).Value = 1.036943704306534
$ws.Cells.Item(10, 5).Value = 1.027328081870708
$ws.Cells.Item(10, 6).Value = 1.044958252949471
$ws.Cells.Item(10, 9).Value = 1.034380707133342
$ws.Cells.Item(10, 10).Value = 1.033520814853591
$ws.Cells.Item(10, 11).Value = 1.04038593452166
$ws.Cells.Item(10, 12).Value = 1.030804953030839
$ws.Cells.Item(10, 13).Value = 1.048372147136781

# Row 11
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.026478564324758
$ws.Cells.Item(11, 4).Value = 1.036443440318482
$ws.Cells.Item(11, 5).Value = 1.026778499453484
$ws.Cells.Item(11, 6).Value = 1.044370928184025
$ws.Cells.Item(11, 9).Value = 1.03426255175727
$ws.Cells.Item(11, 10).Value = 1.033108758598985
$ws.Cells.Item(11, 11).Value = 1.04000890074484
$ws.Cells.Item(11, 12).Value = 1.030380128702086
$ws.Cells.Item(11, 13).Value = 1.047907274216037

# Row 12
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.026237767138234
$ws.Cells.Item(12, 4).Value = 1.036257678137218
$ws.Cells.Item(12, 5).Value = 1.026574520515102
$ws.Cells.Item(12, 6).Value = 1.044152887763824
$ws.Cells.Item(12, 9).Value = 1.034218444024338
$ws.Cells.Item(12, 10).Value = 1.03295569803995
$ws.Cells.Item(12, 11).Value = 1.039868788100282
$ws.Cells.Item(12, 12).Value = 1.030222370356616
$ws.Cells.Item(12, 13).Value = 1.047734599740375

# Row 13
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.026289414257256
$ws.Cells.Item(13, 4).Value = 1.036297522103177
$ws.Cells.Item(13, 5).Value = 1.026618267409102
$ws.Cells.Item(13, 6).Value = 1.044199652805379
$ws.Cells.Item(13, 9).Value = 1.034227915210474
$ws.Cells.Item(13, 10).Value = 1.032988530239534
$ws.Cells.Item(13, 11).Value = 1.039898845686638
$ws.Cells.Item(13, 12).Value = 1.030256208216829
$ws.Cells.Item(13, 13).Value = 1.047771638981722

# Row 14
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.026458657838853
$ws.Cells.Item(14, 4).Value = 1.036428083964837
$ws.Cells.Item(14, 5).Value = 1.026761635210503
$ws.Cells.Item(14, 6).Value = 1.04435290247666
$ws.Cells.Item(14, 9).Value = 1.034258910274455
$ws.Cells.Item(14, 10).Value = 1.033096106652428
$ws.Cells.Item(14, 11).Value = 1.039997320320567
$ws.Cells.Item(14, 12).Value = 1.030367087519789
$ws.Cells.Item(14, 13).Value = 1.047893000882761

# Row 15
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.026562948129963
$ws.Cells.Item(15, 4).Value = 1.036508535181067
$ws.Cells.Item(15, 5).Value = 1.026849990142725
$ws.Cells.Item(15, 6).Value = 1.044447340354313
$ws.Cells.Item(15, 9).Value = 1.034277978276887
$ws.Cells.Item(15, 10).Value = 1.033162387458886
$ws.Cells.Item(15, 11).Value = 1.040057985141569
$ws.Cells.Item(15, 12).Value = 1.030435409290407
$ws.Cells.Item(15, 13).Value = 1.047967775993636

# Row 16
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.027170171834308
$ws.Cells.Item(16, 4).Value = 1.036976913258947
$ws.Cells.Item(16, 5).Value = 1.027364578275709
$ws.Cells.Item(16, 6).Value = 1.044997248252351
$ws.Cells.Item(16, 9).Value = 1.034388517897906
$ws.Cells.Item(16, 10).Value = 1.033548160987128
$ws.Cells.Item(16, 11).Value = 1.040410947769486
$ws.Cells.Item(16, 12).Value = 1.030833152824195
$ws.Cells.Item(16, 13).Value = 1.048402999161248

# Row 17
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.027551250799699
$ws.Cells.Item(17, 4).Value = 1.037270816606136
$ws.Cells.Item(17, 5).Value = 1.027687650360092
$ws.Cells.Item(17, 6).Value = 1.045342399998525
$ws.Cells.Item(17, 9).Value = 1.034457464669391
$ws.Cells.Item(17, 10).Value = 1.03379013751062
$ws.Cells.Item(17, 11).Value = 1.040632234077622
$ws.Cells.Item(17, 12).Value = 1.031082717650552
$ws.Cells.Item(17, 13).Value = 1.048676001949621

# Row 18
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.027773593627206
$ws.Cells.Item(18, 4).Value = 1.037442281837406
$ws.Cells.Item(18, 5).Value = 1.027876194978368
$ws.Cells.Item(18, 6).Value = 1.04554379606542
$ws.Cells.Item(18, 9).Value = 1.034497538789297
$ws.Cells.Item(18, 10).Value = 1.033931274784417
$ws.Cells.Item(18, 11).Value = 1.040761263853301
$ws.Cells.Item(18, 12).Value = 1.03122830980039
$ws.Cells.Item(18, 13).Value = 1.048835239048663

# Row 19
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.027849417942646
$ws.Cells.Item(19, 4).Value = 1.037500753180103
$ws.Cells.Item(19, 5).Value = 1.027940501106426
$ws.Cells.Item(19, 6).Value = 1.045612479580635
$ws.Cells.Item(19, 9).Value = 1.034511179035249
$ws.Cells.Item(19, 10).Value = 1.033979398320029
$ws.Cells.Item(19, 11).Value = 1.040805252384993
$ws.Cells.Item(19, 12).Value = 1.031277957207392
$ws.Cells.Item(19, 13).Value = 1.048889534611862

# Row 20
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.02751035780983
$ws.Cells.Item(20, 4).Value = 1.037239279785612
$ws.Cells.Item(20, 5).Value = 1.027652977209221
$ws.Cells.Item(20, 6).Value = 1.045305360726084
$ws.Cells.Item(20, 9).Value = 1.03445008195176
$ws.Cells.Item(20, 10).Value = 1.033764176079189
$ws.Cells.Item(20, 11).Value = 1.04060849656908
$ws.Cells.Item(20, 12).Value = 1.03105593909203
$ws.Cells.Item(20, 13).Value = 1.048646711409313

# Row 21
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.02640881695971
$ws.Cells.Item(21, 4).Value = 1.036389635135029
$ws.Cells.Item(21, 5).Value = 1.02671941252717
$ws.Cells.Item(21, 6).Value = 1.044307770995502
$ws.Cells.Item(21, 9).Value = 1.034249789057672
$ws.Cells.Item(21, 10).Value = 1.033064428197648
$ws.Cells.Item(21, 11).Value = 1.039968323795701
$ws.Cells.Item(21, 12).Value = 1.030334435207678
$ws.Cells.Item(21, 13).Value = 1.047857262817721

# Row 22
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.025716833580775
$ws.Cells.Item(22, 4).Value = 1.035855766857162
$ws.Cells.Item(22, 5).Value = 1.02613337225249
$ws.Cells.Item(22, 6).Value = 1.043681230994538
$ws.Cells.Item(22, 9).Value = 1.03412258687833
$ws.Cells.Item(22, 10).Value = 1.032624443229522
$ws.Cells.Item(22, 11).Value = 1.03956544322682
$ws.Cells.Item(22, 12).Value = 1.029881030871608
$ws.Cells.Item(22, 13).Value = 1.047360905879432

# Row 23
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.026083610021508
$ws.Cells.Item(23, 4).Value = 1.036138748135349
$ws.Cells.Item(23, 5).Value = 1.026443954817898
$ws.Cells.Item(23, 6).Value = 1.044013306380636
$ws.Cells.Item(23, 9).Value = 1.034190139381379
$ws.Cells.Item(23, 10).Value = 1.032857689722368
$ws.Cells.Item(23, 11).Value = 1.039779053385268
$ws.Cells.Item(23, 12).Value = 1.030121366603583
$ws.Cells.Item(23, 13).Value = 1.047624033626546

# Row 24
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.027528835388089
$ws.Cells.Item(24, 4).Value = 1.037253529805041
$ws.Cells.Item(24, 5).Value = 1.027668644198148
$ws.Cells.Item(24, 6).Value = 1.045322096947688
$ws.Cells.Item(24, 9).Value = 1.03445341832091
$ws.Cells.Item(24, 10).Value = 1.033775906943559
$ws.Cells.Item(24, 11).Value = 1.040619222659702
$ws.Cells.Item(24, 12).Value = 1.031068039092149
$ws.Cells.Item(24, 13).Value = 1.048659946546283

# Row 25
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.029208521893179
$ws.Cells.Item(25, 4).Value = 1.038548594141008
$ws.Cells.Item(25, 5).Value = 1.029093830084577
$ws.Cells.Item(25, 6).Value = 1.046843808403459
$ws.Cells.Item(25, 9).Value = 1.034753391558862
$ws.Cells.Item(25, 10).Value = 1.034841304845953
$ws.Cells.Item(25, 11).Value = 1.041592510355073
$ws.Cells.Item(25, 12).Value = 1.032167585182033
$ws.Cells.Item(25, 13).Value = 1.04986202987562
